$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.990.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.994.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.55%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.16%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.61'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.12%  '

# Row 9
$ws.Range("E9").Value = '  -2.40%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '

# Row 11
$ws.Range("E11").Value = '  -3.77%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0979'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.64%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.44%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.287.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.758'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.82%  '

# Row 17
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.84%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.006.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.912.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0811'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.55%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.01%  '

# Row 30
$ws.Range("E30").Value = '  -2.10%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.59%  '

# Row 32
$ws.Range("E32").Value = '  -1.48%  '

# Row 33
$ws.Range("E33").Value = '  -4.91%  '

# Row 34
$ws.Range("E34").Value = '  -7.01%  '

# Row 35
$ws.Range("E35").Value = '  -5.62%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.66%  '

# Row 37
$ws.Range("E37").Value = '  +0.23%  '

# Row 38
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.68%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.43%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.60%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.431.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.50%  '

# Row 43
$ws.Range("E43").Value = '  -4.01%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0204'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.50%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0891'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.26%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.82%  '

# Row 48
$ws.Range("E48").Value = '  -2.46%  '

# Row 49
$ws.Range("E49").Value = '  +0.44%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.179.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.42%  '
